$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70 timestamp gets a (sub-millisecond) precision update.
$ws.Range("A70").Value = 44414.39481285879

# Append a new log entry as row 71 (mirrors the structure of row 70).
$ws.Range("A71").Value = 44416.82094497623
$ws.Range("A71").NumberFormat = $ws.Range("A70").NumberFormat

$ws.Range("B71").Value = "work"
$ws.Range("C71").Value = 1
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 1
$ws.Range("F71").Value = "other2"
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4.09
$ws.Range("K71").Value = "Tier Summer Pass"
